# Auto-generated script applying scheduled-runner market-data refresh
# to the Leve profit tables (columns H-N) across all 8 job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Cells.Item(7, 8).Value = 25806.5
$ws.Cells.Item(7, 10).Value = 26108.666
$ws.Cells.Item(7, 12).Value = 26108.666
$ws.Cells.Item(7, 14).Value = -26332.666
$ws.Cells.Item(14, 8).Value = 25806.5
$ws.Cells.Item(14, 10).Value = 26108.666
$ws.Cells.Item(14, 12).Value = 26108.666
$ws.Cells.Item(14, 14).Value = -26490.666
$ws.Cells.Item(38, 8).Value = 3504.818
$ws.Cells.Item(38, 9).Value = 129.83333
$ws.Cells.Item(38, 10).Value = 7554.8
$ws.Cells.Item(38, 11).Value = 389.49999
$ws.Cells.Item(38, 12).Value = 22664.4
$ws.Cells.Item(38, 13).Value = -17.49998999999997
$ws.Cells.Item(38, 14).Value = -23408.4
$ws.Cells.Item(58, 8).Value = 101.333336
$ws.Cells.Item(58, 9).Value = 101.333336
$ws.Cells.Item(58, 11).Value = 304.000008
$ws.Cells.Item(58, 13).Value = -154.000008
$ws.Cells.Item(64, 8).Value = 100002940
$ws.Cells.Item(64, 10).Value = 3616.5
$ws.Cells.Item(64, 12).Value = 3616.5
$ws.Cells.Item(64, 14).Value = -4112.5
$ws.Cells.Item(67, 8).Value = 100002940
$ws.Cells.Item(67, 10).Value = 3616.5
$ws.Cells.Item(67, 12).Value = 3616.5
$ws.Cells.Item(67, 14).Value = -5332.5
$ws.Cells.Item(74, 8).Value = 38465732
$ws.Cells.Item(74, 9).Value = 50003652
$ws.Cells.Item(74, 10).Value = 6000
$ws.Cells.Item(74, 11).Value = 50003652
$ws.Cells.Item(74, 12).Value = 6000
$ws.Cells.Item(74, 13).Value = -50002716
$ws.Cells.Item(74, 14).Value = -7872
$ws.Cells.Item(77, 8).Value = 38465732
$ws.Cells.Item(77, 9).Value = 50003652
$ws.Cells.Item(77, 10).Value = 6000
$ws.Cells.Item(77, 11).Value = 250018260
$ws.Cells.Item(77, 12).Value = 30000
$ws.Cells.Item(77, 13).Value = -250013580
$ws.Cells.Item(77, 14).Value = -39360
$ws.Cells.Item(82, 8).Value = 2173.1667
$ws.Cells.Item(82, 9).Value = 2173.1667
$ws.Cells.Item(82, 11).Value = 6519.500100000001
$ws.Cells.Item(82, 13).Value = -6113.500100000001
$ws.Cells.Item(85, 8).Value = 2173.1667
$ws.Cells.Item(85, 9).Value = 2173.1667
$ws.Cells.Item(85, 11).Value = 6519.500100000001
$ws.Cells.Item(85, 13).Value = -5115.500100000001
$ws.Cells.Item(106, 8).Value = 2754.7334
$ws.Cells.Item(106, 9).Value = 3447.8572
$ws.Cells.Item(106, 10).Value = 2148.25
$ws.Cells.Item(106, 11).Value = 3447.8572
$ws.Cells.Item(106, 12).Value = 2148.25
$ws.Cells.Item(106, 13).Value = -2816.8572
$ws.Cells.Item(106, 14).Value = -3410.25
$ws.Cells.Item(111, 8).Value = 5809.5
$ws.Cells.Item(111, 9).Value = 3572
$ws.Cells.Item(111, 10).Value = 16997
$ws.Cells.Item(111, 11).Value = 10716
$ws.Cells.Item(111, 12).Value = 50991
$ws.Cells.Item(111, 13).Value = -7649
$ws.Cells.Item(111, 14).Value = -57125
$ws.Cells.Item(132, 8).Value = 5099.523
$ws.Cells.Item(132, 9).Value = 1425.9474
$ws.Cells.Item(132, 11).Value = 4277.8422
$ws.Cells.Item(132, 13).Value = -1747.8422
$ws.Cells.Item(135, 8).Value = 3033.6538
$ws.Cells.Item(135, 9).Value = 790.13635
$ws.Cells.Item(135, 11).Value = 7111.22715
$ws.Cells.Item(135, 13).Value = -4576.22715
$ws.Cells.Item(137, 8).Value = 231364.45
$ws.Cells.Item(137, 9).Value = 374253.22
$ws.Cells.Item(137, 10).Value = 4423.4707
$ws.Cells.Item(137, 11).Value = 1122759.66
$ws.Cells.Item(137, 12).Value = 13270.4121
$ws.Cells.Item(137, 13).Value = -1120209.66
$ws.Cells.Item(137, 14).Value = -18370.4121

$ws = $wb.Sheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 606
$ws.Cells.Item(74, 9).Value = 431.8125
$ws.Cells.Item(74, 11).Value = 431.8125
$ws.Cells.Item(74, 13).Value = 442.1875
$ws.Cells.Item(77, 8).Value = 606
$ws.Cells.Item(77, 9).Value = 431.8125
$ws.Cells.Item(77, 11).Value = 2159.0625
$ws.Cells.Item(77, 13).Value = 2208.9375
$ws.Cells.Item(97, 8).Value = 15656.65
$ws.Cells.Item(97, 9).Value = 20341.066
$ws.Cells.Item(97, 10).Value = 1603.4
$ws.Cells.Item(97, 11).Value = 20341.066
$ws.Cells.Item(97, 12).Value = 1603.4
$ws.Cells.Item(97, 13).Value = -19845.066
$ws.Cells.Item(97, 14).Value = -2595.4
$ws.Cells.Item(102, 8).Value = 1396.8096
$ws.Cells.Item(102, 9).Value = 1346.75
$ws.Cells.Item(102, 11).Value = 1346.75
$ws.Cells.Item(102, 13).Value = 275.25
$ws.Cells.Item(110, 8).Value = 1674.5
$ws.Cells.Item(110, 9).Value = 1674.5
$ws.Cells.Item(110, 11).Value = 1674.5
$ws.Cells.Item(110, 13).Value = 370.5
$ws.Cells.Item(122, 8).Value = 3246.2415
$ws.Cells.Item(122, 9).Value = 2045.64
$ws.Cells.Item(122, 11).Value = 6136.92
$ws.Cells.Item(122, 13).Value = -3686.92

$ws = $wb.Sheets.Item("BSM")
$ws.Cells.Item(64, 8).Value = 712.25
$ws.Cells.Item(64, 9).Value = 91.5
$ws.Cells.Item(64, 10).Value = 1333
$ws.Cells.Item(64, 11).Value = 91.5
$ws.Cells.Item(64, 12).Value = 1333
$ws.Cells.Item(64, 13).Value = 133.5
$ws.Cells.Item(64, 14).Value = -1783
$ws.Cells.Item(67, 8).Value = 712.25
$ws.Cells.Item(67, 9).Value = 91.5
$ws.Cells.Item(67, 10).Value = 1333
$ws.Cells.Item(67, 11).Value = 91.5
$ws.Cells.Item(67, 12).Value = 1333
$ws.Cells.Item(67, 13).Value = 688.5
$ws.Cells.Item(67, 14).Value = -2893
$ws.Cells.Item(86, 8).Value = 3057.0435
$ws.Cells.Item(86, 9).Value = 2210.5293
$ws.Cells.Item(86, 10).Value = 5455.5
$ws.Cells.Item(86, 11).Value = 2210.5293
$ws.Cells.Item(86, 12).Value = 5455.5
$ws.Cells.Item(86, 13).Value = -1087.5293
$ws.Cells.Item(86, 14).Value = -7701.5
$ws.Cells.Item(89, 8).Value = 3057.0435
$ws.Cells.Item(89, 9).Value = 2210.5293
$ws.Cells.Item(89, 10).Value = 5455.5
$ws.Cells.Item(89, 11).Value = 11052.6465
$ws.Cells.Item(89, 12).Value = 27277.5
$ws.Cells.Item(89, 13).Value = -5436.646500000001
$ws.Cells.Item(89, 14).Value = -38509.5
$ws.Cells.Item(134, 8).Value = 5131.5093
$ws.Cells.Item(134, 9).Value = 2100.6538
$ws.Cells.Item(134, 11).Value = 6301.9614
$ws.Cells.Item(134, 13).Value = -3766.9614

$ws = $wb.Sheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 2509.5557
$ws.Cells.Item(16, 9).Value = 1581.8572
$ws.Cells.Item(16, 11).Value = 1581.8572
$ws.Cells.Item(16, 13).Value = -1294.8572
$ws.Cells.Item(62, 8).Value = 52987.5
$ws.Cells.Item(62, 9).Value = 25650
$ws.Cells.Item(62, 10).Value = 135000
$ws.Cells.Item(62, 11).Value = 25650
$ws.Cells.Item(62, 12).Value = 135000
$ws.Cells.Item(62, 13).Value = -25026
$ws.Cells.Item(62, 14).Value = -136248
$ws.Cells.Item(65, 8).Value = 52987.5
$ws.Cells.Item(65, 9).Value = 25650
$ws.Cells.Item(65, 10).Value = 135000
$ws.Cells.Item(65, 11).Value = 128250
$ws.Cells.Item(65, 12).Value = 675000
$ws.Cells.Item(65, 13).Value = -125130
$ws.Cells.Item(65, 14).Value = -681240
$ws.Cells.Item(113, 8).Value = 2509.5557
$ws.Cells.Item(113, 9).Value = 1581.8572
$ws.Cells.Item(113, 11).Value = 1581.8572
$ws.Cells.Item(113, 13).Value = 588.1428000000001
$ws.Cells.Item(141, 8).Value = 180624.9
$ws.Cells.Item(141, 10).Value = 189694.39
$ws.Cells.Item(141, 12).Value = 189694.39
$ws.Cells.Item(141, 14).Value = -200054.39

$ws = $wb.Sheets.Item("CUL")
$ws.Cells.Item(19, 8).Value = 0
$ws.Cells.Item(19, 10).Value = 0
$ws.Cells.Item(19, 12).ClearContents()
$ws.Cells.Item(19, 14).Value = 0
$ws.Cells.Item(114, 8).Value = 2013.875
$ws.Cells.Item(114, 9).Value = 363
$ws.Cells.Item(114, 10).Value = 2249.7144
$ws.Cells.Item(114, 11).Value = 1089
$ws.Cells.Item(114, 12).Value = 6749.1432
$ws.Cells.Item(114, 13).Value = 2165
$ws.Cells.Item(114, 14).Value = -13257.1432
$ws.Cells.Item(131, 8).Value = 1492629.9
$ws.Cells.Item(131, 9).Value = 1035.4546
$ws.Cells.Item(131, 11).Value = 3106.3638
$ws.Cells.Item(131, 13).Value = 1933.6362
$ws.Cells.Item(140, 8).Value = 10399.667
$ws.Cells.Item(140, 9).Value = 13999.5
$ws.Cells.Item(140, 10).Value = 3200
$ws.Cells.Item(140, 11).Value = 41998.5
$ws.Cells.Item(140, 12).Value = 9600
$ws.Cells.Item(140, 13).Value = -36818.5
$ws.Cells.Item(140, 14).Value = -19960

$ws = $wb.Sheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 217471.02
$ws.Cells.Item(122, 9).Value = 360704.3
$ws.Cells.Item(122, 11).Value = 1082112.9
$ws.Cells.Item(122, 13).Value = -1079662.9
$ws.Cells.Item(123, 8).Value = 40072.6
$ws.Cells.Item(123, 10).Value = 40072.6
$ws.Cells.Item(123, 12).Value = 40072.6
$ws.Cells.Item(123, 14).Value = -44972.6

$ws = $wb.Sheets.Item("LTW")
$ws.Cells.Item(19, 8).Value = 61910
$ws.Cells.Item(19, 9).Value = 61910
$ws.Cells.Item(19, 10).Value = 0
$ws.Cells.Item(19, 11).Value = 61910
$ws.Cells.Item(19, 12).Value = 0
$ws.Cells.Item(19, 13).ClearContents()
$ws.Cells.Item(19, 14).Value = -61740
$ws.Cells.Item(22, 8).Value = 988
$ws.Cells.Item(22, 9).Value = 778
$ws.Cells.Item(22, 10).Value = 1450
$ws.Cells.Item(22, 11).Value = 778
$ws.Cells.Item(22, 12).Value = 1450
$ws.Cells.Item(22, 13).Value = -483
$ws.Cells.Item(22, 14).Value = -2040
$ws.Cells.Item(26, 8).Value = 34999.332
$ws.Cells.Item(26, 9).Value = 24999
$ws.Cells.Item(26, 10).Value = 39999.5
$ws.Cells.Item(26, 11).Value = 24999
$ws.Cells.Item(26, 12).Value = 39999.5
$ws.Cells.Item(26, 13).Value = -24704
$ws.Cells.Item(26, 14).Value = -40589.5
$ws.Cells.Item(27, 8).Value = 988
$ws.Cells.Item(27, 9).Value = 778
$ws.Cells.Item(27, 10).Value = 1450
$ws.Cells.Item(27, 11).Value = 778
$ws.Cells.Item(27, 12).Value = 1450
$ws.Cells.Item(27, 13).Value = -671
$ws.Cells.Item(27, 14).Value = -1664
$ws.Cells.Item(41, 8).Value = 0
$ws.Cells.Item(41, 10).Value = 0
$ws.Cells.Item(41, 12).ClearContents()
$ws.Cells.Item(41, 14).Value = 0
$ws.Cells.Item(46, 8).Value = 4671.6113
$ws.Cells.Item(46, 10).Value = 4358.1763
$ws.Cells.Item(46, 12).Value = 4358.1763
$ws.Cells.Item(46, 14).Value = -4734.1763
$ws.Cells.Item(47, 8).Value = 36021.332
$ws.Cells.Item(47, 9).Value = 19999
$ws.Cells.Item(47, 10).Value = 44032.5
$ws.Cells.Item(47, 11).Value = 19999
$ws.Cells.Item(47, 12).Value = 44032.5
$ws.Cells.Item(47, 13).Value = -19509
$ws.Cells.Item(47, 14).Value = -45012.5
$ws.Cells.Item(48, 8).Value = 0
$ws.Cells.Item(48, 9).Value = 0
$ws.Cells.Item(48, 11).Value = 0
$ws.Cells.Item(48, 13).ClearContents()
$ws.Cells.Item(52, 8).Value = 36021.332
$ws.Cells.Item(52, 9).Value = 19999
$ws.Cells.Item(52, 10).Value = 44032.5
$ws.Cells.Item(52, 11).Value = 19999
$ws.Cells.Item(52, 12).Value = 44032.5
$ws.Cells.Item(52, 13).Value = -19766
$ws.Cells.Item(52, 14).Value = -44498.5
$ws.Cells.Item(55, 8).Value = 550.8
$ws.Cells.Item(55, 9).Value = 182
$ws.Cells.Item(55, 11).Value = 182
$ws.Cells.Item(55, 13).Value = -9
$ws.Cells.Item(98, 8).Value = 50000
$ws.Cells.Item(98, 10).Value = 50000
$ws.Cells.Item(98, 12).Value = 50000
$ws.Cells.Item(98, 14).Value = -55990

$ws = $wb.Sheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 2551.8
$ws.Cells.Item(132, 9).Value = 939.9259
$ws.Cells.Item(132, 11).Value = 2819.7777
$ws.Cells.Item(132, 13).Value = -289.7776999999996

